$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.264496
$ws.Range("H2").Value = 0.793488
$ws.Range("I2").Value = 0.001006353962629067
$ws.Range("J2").Value = 0.001006353962629067
$ws.Range("M2").Value = 0.06166766666666667
$ws.Range("N2").Value = 0.185003
$ws.Range("O2").Value = 0.3189772891852935
$ws.Range("P2").Value = 0.3189772891852935
$ws.Range("Q2").Value = 0.01631085116266667
$ws.Range("R2").Value = 0.146797660464
$ws.Range("S2").Value = 0.000321004058960298
$ws.Range("T2").Value = 0.0003210040589602979

# Row 3
$ws.Range("G3").Value = 0.264496
$ws.Range("H3").Value = 0.793488
$ws.Range("I3").Value = 0.001006353962629067
$ws.Range("J3").Value = 0.001006353962629067
$ws.Range("O3").Value = 0.4045463009579509
$ws.Range("P3").Value = 0.4045463009579509
$ws.Range("Q3").Value = 0.02068640849066667
$ws.Range("R3").Value = 0.186177676416
$ws.Range("S3").Value = 0.000407116773035965
$ws.Range("T3").Value = 0.000407116773035965

# Row 4
$ws.Range("G4").Value = 0.264496
$ws.Range("H4").Value = 0.793488
$ws.Range("I4").Value = 0.001006353962629067
$ws.Range("J4").Value = 0.001006353962629067
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05345100000000001
$ws.Range("N4").Value = 0.160353
$ws.Range("O4").Value = 0.2764764098567557
$ws.Range("P4").Value = 0.2764764098567557
$ws.Range("Q4").Value = 0.014137575696
$ws.Range("R4").Value = 0.127238181264
$ws.Range("S4").Value = 0.0002782331306328042
$ws.Range("T4").Value = 0.0002782331306328041

# Row 5
$ws.Range("H5").Value = 737.537796
$ws.Range("I5").Value = 0.9353942133886188
$ws.Range("J5").Value = 0.935394213388619
$ws.Range("M5").Value = 0.06166766666666667
$ws.Range("N5").Value = 0.185003
$ws.Range("O5").Value = 0.3189772891852935
$ws.Range("P5").Value = 0.3189772891852935
$ws.Range("Q5").Value = 15.160744985932
$ws.Range("R5").Value = 136.446704873388
$ws.Range("S5").Value = 0.2983695105063116
$ws.Range("T5").Value = 0.2983695105063117

# Row 6
$ws.Range("H6").Value = 737.537796
$ws.Range("I6").Value = 0.9353942133886188
$ws.Range("J6").Value = 0.935394213388619
$ws.Range("O6").Value = 0.4045463009579509
$ws.Range("P6").Value = 0.4045463009579509
$ws.Range("S6").Value = 0.3784102689638379
$ws.Range("T6").Value = 0.378410268963838

# Row 7
$ws.Range("H7").Value = 737.537796
$ws.Range("I7").Value = 0.9353942133886188
$ws.Range("J7").Value = 0.935394213388619
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.05345100000000001
$ws.Range("N7").Value = 0.160353
$ws.Range("O7").Value = 0.2764764098567557
$ws.Range("P7").Value = 0.2764764098567557
$ws.Range("Q7").Value = 13.140710911332
$ws.Range("R7").Value = 118.266398201988
$ws.Range("S7").Value = 0.2586144339184694
$ws.Range("T7").Value = 0.2586144339184694

# Row 8
$ws.Range("G8").Value = 16.71558533333333
$ws.Range("H8").Value = 50.146756
$ws.Range("I8").Value = 0.06359943264875202
$ws.Range("J8").Value = 0.06359943264875202
$ws.Range("M8").Value = 0.06166766666666667
$ws.Range("N8").Value = 0.185003
$ws.Range("O8").Value = 0.3189772891852935
$ws.Range("P8").Value = 0.3189772891852935
$ws.Range("Q8").Value = 1.030811144474222
$ws.Range("R8").Value = 9.277300300268
$ws.Range("S8").Value = 0.02028677462002157
$ws.Range("T8").Value = 0.02028677462002157

# Row 9
$ws.Range("G9").Value = 16.71558533333333
$ws.Range("H9").Value = 50.146756
$ws.Range("I9").Value = 0.06359943264875202
$ws.Range("J9").Value = 0.06359943264875202
$ws.Range("O9").Value = 0.4045463009579509
$ws.Range("P9").Value = 0.4045463009579509
$ws.Range("Q9").Value = 1.307337072643556
$ws.Range("R9").Value = 11.766033653792
$ws.Range("S9").Value = 0.02572891522107696
$ws.Range("T9").Value = 0.02572891522107696

# Row 10
$ws.Range("G10").Value = 16.71558533333333
$ws.Range("H10").Value = 50.146756
$ws.Range("I10").Value = 0.06359943264875202
$ws.Range("J10").Value = 0.06359943264875202
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.05345100000000001
$ws.Range("N10").Value = 0.160353
$ws.Range("O10").Value = 0.2764764098567557
$ws.Range("P10").Value = 0.2764764098567557
$ws.Range("Q10").Value = 0.8934647516520001
$ws.Range("R10").Value = 8.041182764868001
$ws.Range("S10").Value = 0.01758374280765349
$ws.Range("T10").Value = 0.01758374280765349

